$wb = $excel.ActiveWorkbook

# --- Original sheets ---
# Sheet1 (currently active) holds the Journal ledger data -> rename to "Journal"
# Sheet2 is an empty placeholder sheet -> remove it, we'll add fresh report sheets instead
$wb.Worksheets.Item(1).Name = "Journal"

$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

# --- Add the four new (still blank/template) report sheets, all positioned
#     before the Journal sheet. Each Add(Before) call inserts the new sheet
#     immediately in front of whatever "Journal" resolves to at that moment,
#     so re-resolving "Journal" by name (sheet references are live/index-bound
#     and go stale once the collection is reshuffled) and adding in the same
#     left-to-right order as the desired result keeps everything in place:
#     Statement of Owner Equity, Income Statement, Balance Sheet, T Accounts, Journal ---
$ownerEquity = $wb.Worksheets.Add($wb.Worksheets.Item("Journal"))
$ownerEquity.Name = "Statement of Owner Equity"

$incomeStatement = $wb.Worksheets.Add($wb.Worksheets.Item("Journal"))
$incomeStatement.Name = "Income Statement"

$balanceSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Journal"))
$balanceSheet.Name = "Balance Sheet"

$tAccounts = $wb.Worksheets.Add($wb.Worksheets.Item("Journal"))
$tAccounts.Name = "T Accounts"

# --- Populate the header row of each new report sheet with the same
#     placeholder header used in the source export (ID / FILL IN / FILL IN) ---
foreach ($sheetName in @("Statement of Owner Equity", "Income Statement", "Balance Sheet", "T Accounts")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range("A1").Value = "ID"
    $sheet.Range("B1").Value = "FILL IN"
    $sheet.Range("C1").Value = "FILL IN"
}

# --- Final tab order is: Statement of Owner Equity, Income Statement,
#     Balance Sheet, T Accounts, Journal -- with Journal left as the active tab ---
$wb.Worksheets.Item("Journal").Activate()
